$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values for the three remaining "MuSCs" sending-cluster rows.
# Row 2: MuSCs -> Cdh1/Itgb7 -> ECs
# Row 3: MuSCs -> Cdh1/Itgb7 -> FAPs
# Row 4: MuSCs -> Cdh1/Itgb7 -> MuSCs
$data = @(
    @("MuSCs","Cdh1","Itgb7","ECs",3,1,0.9477166666666667,2.84315,1,1,3,1,0.8194946666666666,2.458484,0.1466535424263973,0.1466535424263973,0.7766487538444443,6.9898387846,0.1466535424263973,0.1466535424263973),
    @("MuSCs","Cdh1","Itgb7","FAPs",3,1,0.9477166666666667,2.84315,1,1,3,1,3.605649999999999,10.81695,0.6452529427684778,0.6452529427684778,3.417134599166666,30.7542113925,0.6452529427684778,0.6452529427684778),
    @("MuSCs","Cdh1","Itgb7","MuSCs",3,1,0.9477166666666667,2.84315,1,1,3,1,1.162819,3.488457,0.2080935148051249,0.2080935148051249,1.102022946616667,9.918206519550001,0.2080935148051249,0.2080935148051249)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $rowVals = $data[$i]
    for ($c = 1; $c -le $rowVals.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}

# Remove the now-unused trailing rows (old rows 5-7).
$ws.Range("A5:T7").Delete() | Out-Null
